$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (KILLS) and E (DEATHS) were stored as text-that-looks-like-a-number
# (and in a couple of spots outright typos: "a" and "erro"). This pass fixes
# the running kill/death tally so every cell holds a real number.
$aValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 3
    22 = 3
    23 = 4
    24 = 4
    25 = 4
    26 = 4
    27 = 4
    28 = 5
    29 = 5
    30 = 5
    31 = 5
    32 = 5
    33 = 5
    34 = 5
    35 = 5
    36 = 5
    37 = 5
    38 = 5
    39 = 5
    40 = 5
    41 = 5
}

$eValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 2
    25 = 3
    26 = 3
    27 = 3
    28 = 3
    29 = 3
    30 = 3
    31 = 3
    32 = 3
    33 = 3
    34 = 3
    35 = 3
    36 = 3
    37 = 4
    38 = 4
    39 = 4
    40 = 4
    41 = 4
}

foreach ($r in 2..41) {
    $ws.Cells.Item($r, 1).Value = $aValues[$r]
    $ws.Cells.Item($r, 5).Value = $eValues[$r]
}

# Column F (ASSISTS) row 41 had the literal typo "erro" — replace with the
# correct running-assist count "3". It must stay a TEXT cell (like its
# neighbours in column F), so force text entry via a temporary Text format,
# then restore the cell to the default (unstyled) look.
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "3"
$ws.Range("F41").Style = "Normal"
